$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A21").Value = "Naujų skelbimų garsas"
$ws.Range("B21").Value = "Ignas"
$ws.Range("C21").Value = "vidutinis"
$ws.Range("C21").HorizontalAlignment = -4108

$ws.Range("C27").Select()
